$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 13:42:17.449699",
    "2021-10-05 13:42:17.449713",
    "2021-10-05 13:42:17.449717",
    "2021-10-05 13:42:17.449720",
    "2021-10-05 13:42:17.449723",
    "2021-10-05 13:42:17.449727",
    "2021-10-05 13:42:17.449730",
    "2021-10-05 13:42:17.449733",
    "2021-10-05 13:42:17.449736",
    "2021-10-05 13:42:17.449740",
    "2021-10-05 13:42:17.449742",
    "2021-10-05 13:42:17.449746",
    "2021-10-05 13:42:17.449749",
    "2021-10-05 13:42:17.449751",
    "2021-10-05 13:42:17.449754",
    "2021-10-05 13:42:17.449757",
    "2021-10-05 13:42:17.449761",
    "2021-10-05 13:42:17.449764",
    "2021-10-05 13:42:17.449767",
    "2021-10-05 13:42:17.449770",
    "2021-10-05 13:42:17.449773",
    "2021-10-05 13:42:17.449776"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
